# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.073.40"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.32%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.73%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.017"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.94%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06439"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07781"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.644.85"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.260"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.864.29"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5469"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7959"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.050.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.017"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "205.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.319"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.982"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.017"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.978"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +14.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.26%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1156"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.58%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.820"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05046"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.246"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.275"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.82%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.212"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.548"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8950"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.40%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.618"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.16%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5658"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.126.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01566"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.590"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.91%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.019"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.653"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8200"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.97"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.773.14"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4564"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.015"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05053"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.25%  "
